$d = $word.ActiveDocument

$d.Content.Find.Execute("110÷2=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "315÷5=63, 0", 2) | Out-Null
$d.Content.Find.Execute("731÷4=182, 3", $true, $false, $false, $false, $false, $true, 1, $false, "342÷7=48, 6", 2) | Out-Null
$d.Content.Find.Execute("219÷9=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "860÷8=107, 4", 2) | Out-Null
$d.Content.Find.Execute("759÷7=108, 3", $true, $false, $false, $false, $false, $true, 1, $false, "622÷5=124, 2", 2) | Out-Null
$d.Content.Find.Execute("460÷7=65, 5", $true, $false, $false, $false, $false, $true, 1, $false, "865÷2=432, 1", 2) | Out-Null
$d.Content.Find.Execute("461÷8=57, 5", $true, $false, $false, $false, $false, $true, 1, $false, "902÷9=100, 2", 2) | Out-Null
$d.Content.Find.Execute("863÷4=215, 3", $true, $false, $false, $false, $false, $true, 1, $false, "302÷8=37, 6", 2) | Out-Null
$d.Content.Find.Execute("613÷8=76, 5", $true, $false, $false, $false, $false, $true, 1, $false, "614÷9=68, 2", 2) | Out-Null
$d.Content.Find.Execute("575÷7=82, 1", $true, $false, $false, $false, $false, $true, 1, $false, "131÷5=26, 1", 2) | Out-Null
$d.Content.Find.Execute("326÷6=54, 2", $true, $false, $false, $false, $false, $true, 1, $false, "739÷2=369, 1", 2) | Out-Null
$d.Content.Find.Execute("748÷7=106, 6", $true, $false, $false, $false, $false, $true, 1, $false, "869÷7=124, 1", 2) | Out-Null
$d.Content.Find.Execute("962÷6=160, 2", $true, $false, $false, $false, $false, $true, 1, $false, "723÷4=180, 3", 2) | Out-Null
$d.Content.Find.Execute("646÷5=129, 1", $true, $false, $false, $false, $false, $true, 1, $false, "557÷6=92, 5", 2) | Out-Null
$d.Content.Find.Execute("673÷9=74, 7", $true, $false, $false, $false, $false, $true, 1, $false, "325÷9=36, 1", 2) | Out-Null
$d.Content.Find.Execute("293÷9=32, 5", $true, $false, $false, $false, $false, $true, 1, $false, "536÷5=107, 1", 2) | Out-Null
$d.Content.Find.Execute("641÷9=71, 2", $true, $false, $false, $false, $false, $true, 1, $false, "811÷3=270, 1", 2) | Out-Null
$d.Content.Find.Execute("872÷4=218, 0", $true, $false, $false, $false, $false, $true, 1, $false, "326÷5=65, 1", 2) | Out-Null
$d.Content.Find.Execute("154÷9=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "376÷7=53, 5", 2) | Out-Null
$d.Content.Find.Execute("340÷2=170, 0", $true, $false, $false, $false, $false, $true, 1, $false, "918÷5=183, 3", 2) | Out-Null
$d.Content.Find.Execute("685÷3=228, 1", $true, $false, $false, $false, $false, $true, 1, $false, "107÷6=17, 5", 2) | Out-Null
$d.Content.Find.Execute("946÷5=189, 1", $true, $false, $false, $false, $false, $true, 1, $false, "698÷3=232, 2", 2) | Out-Null
$d.Content.Find.Execute("960÷5=192, 0", $true, $false, $false, $false, $false, $true, 1, $false, "231÷9=25, 6", 2) | Out-Null
$d.Content.Find.Execute("951÷3=317, 0", $true, $false, $false, $false, $false, $true, 1, $false, "443÷4=110, 3", 2) | Out-Null
$d.Content.Find.Execute("701÷6=116, 5", $true, $false, $false, $false, $false, $true, 1, $false, "805÷8=100, 5", 2) | Out-Null
$d.Content.Find.Execute("397÷9=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "223÷4=55, 3", 2) | Out-Null
